$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: update Hydrogen demand for Iron & steel (B3)
$ws.Range("B3").Value = 3480806.526046263

# Row 3: D3 (Non-metallic minerals / Hydrogen) no longer has a value
$ws.Range("D3").ClearContents()

# Row 4: update Methanol demand for Chemicals (C4)
$ws.Range("C4").Value = 66.31453748371544

# Row 5: update Ammonia demand for Chemicals (C5)
$ws.Range("C5").Value = 1868.224721276874

# Row 7: rename "Other" to "Biogas" and update its value
$ws.Range("A7").Value = "Biogas"
$ws.Range("D7").Value = 2587.885461925465

# Row 8: new "Other" row, formatted like row 7, with a value in D8
$ws.Range("A7:D7").Copy()
$ws.Range("A8:D8").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A8").Value = "Other"
$ws.Range("D8").Value = 3150.187038115711
